# "Austausch Mitarbeiterlaptop" example extended:
# - add "X" marks to existing rows (Incidentmanagement, IT-Kapitalmanagement,
#   Problemmanagement, Releasemanagement)
# - add a new "IT-Assetmanagement" annotation next to IT-Kapitalmanagement
# - add "(X)" marks for Problemmanagement / Releasemanagement
# - append two new activity rows: Deploymentmanagement, Servicevalidierungs- & ~test

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Incidentmanagement
$ws.Range("B8").Value = "X"
$ws.Range("D8").Value = "X"

# Row 9 - IT-Kapitalmanagement
$ws.Range("C9").Value = "X"
$ws.Range("D9").Value = "X"
$ws.Range("F9").Value = "IT-Assetmanagement"

# Row 10 - Problemmanagement
$ws.Range("B10").Value = "(X)"
$ws.Range("D10").Value = "X"

# Row 11 - Releasemanagement
$ws.Range("D11").Value = "X"
$ws.Range("E11").Value = "(X)"

# Row 12 - new activity: Deploymentmanagement
$ws.Range("A12").Value = "Deploymentmanagement"
$ws.Range("D12").Value = "X"

# Row 13 - new activity: Servicevalidierungs- & ~test
$ws.Range("A13").Value = "Servicevalidierungs- & ~test"
$ws.Range("C13").Value = "X"
$ws.Range("D13").Value = "X"
$ws.Range("E13").Value = "(X)"
# Column E has no column-level style, so a brand new row-13 cell needs the
# centered "mark" style (matching E3:E11) applied explicitly.
$ws.Range("E13").HorizontalAlignment = -4108  # xlCenter

# Widen the new annotation column F (~18.68 characters; COM ColumnWidth
# quantizes to whole pixels, so 17.8 is the closest input that lands on it)
$ws.Columns.Item(6).ColumnWidth = 17.8

# Match the author's final selection
$ws.Range("E13").Select()
